$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update time_taken (column F) timestamps for rows 2-70 to match the refreshed panel query time
$timestamps = @(
    "2021-10-05 14:35:10.856506",
    "2021-10-05 14:35:10.856514",
    "2021-10-05 14:35:10.856517",
    "2021-10-05 14:35:10.856520",
    "2021-10-05 14:35:10.856522",
    "2021-10-05 14:35:10.856525",
    "2021-10-05 14:35:10.856527",
    "2021-10-05 14:35:10.856530",
    "2021-10-05 14:35:10.856532",
    "2021-10-05 14:35:10.856535",
    "2021-10-05 14:35:10.856537",
    "2021-10-05 14:35:10.856540",
    "2021-10-05 14:35:10.856542",
    "2021-10-05 14:35:10.856545",
    "2021-10-05 14:35:10.856547",
    "2021-10-05 14:35:10.856550",
    "2021-10-05 14:35:10.856552",
    "2021-10-05 14:35:10.856555",
    "2021-10-05 14:35:10.856557",
    "2021-10-05 14:35:10.856560",
    "2021-10-05 14:35:10.856562",
    "2021-10-05 14:35:10.856565",
    "2021-10-05 14:35:10.856567",
    "2021-10-05 14:35:10.856569",
    "2021-10-05 14:35:10.856572",
    "2021-10-05 14:35:10.856575",
    "2021-10-05 14:35:10.856577",
    "2021-10-05 14:35:10.856579",
    "2021-10-05 14:35:10.856582",
    "2021-10-05 14:35:10.856584",
    "2021-10-05 14:35:10.856587",
    "2021-10-05 14:35:10.856589",
    "2021-10-05 14:35:10.856592",
    "2021-10-05 14:35:10.856595",
    "2021-10-05 14:35:10.856597",
    "2021-10-05 14:35:10.856599",
    "2021-10-05 14:35:10.856602",
    "2021-10-05 14:35:10.856604",
    "2021-10-05 14:35:10.856607",
    "2021-10-05 14:35:10.856609",
    "2021-10-05 14:35:10.856612",
    "2021-10-05 14:35:10.856614",
    "2021-10-05 14:35:10.856617",
    "2021-10-05 14:35:10.856619",
    "2021-10-05 14:35:10.856622",
    "2021-10-05 14:35:10.856624",
    "2021-10-05 14:35:10.856627",
    "2021-10-05 14:35:10.856629",
    "2021-10-05 14:35:10.856632",
    "2021-10-05 14:35:10.856634",
    "2021-10-05 14:35:10.856636",
    "2021-10-05 14:35:10.856639",
    "2021-10-05 14:35:10.856642",
    "2021-10-05 14:35:10.856644",
    "2021-10-05 14:35:10.856647",
    "2021-10-05 14:35:10.856649",
    "2021-10-05 14:35:10.856651",
    "2021-10-05 14:35:10.856654",
    "2021-10-05 14:35:10.856656",
    "2021-10-05 14:35:10.856659",
    "2021-10-05 14:35:10.856661",
    "2021-10-05 14:35:10.856664",
    "2021-10-05 14:35:10.856666",
    "2021-10-05 14:35:10.856668",
    "2021-10-05 14:35:10.856672",
    "2021-10-05 14:35:10.856674",
    "2021-10-05 14:35:10.856677",
    "2021-10-05 14:35:10.856679",
    "2021-10-05 14:35:10.856682"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add a new "metadata" sheet and move it right after the "data" sheet
$meta = $wb.Worksheets.Add()
$meta.Name = "metadata"
$meta.Move($null, $wb.Worksheets.Item("data"))

# Re-fetch live references now that the sheet order has changed
$ws = $wb.Worksheets.Item("data")
$meta = $wb.Worksheets.Item("metadata")

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"
$meta.Range("B1:G1").Style = $ws.Range("B1:B1").Style

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 1).Style = $ws.Range("A2:A2").Style
$meta.Cells.Item(2, 2).Value = "Palmoplantar keratodermas_GEL"
$meta.Cells.Item(2, 3).Value = 3286
$meta.Cells.Item(2, 4).Value = "0.2"
$meta.Cells.Item(2, 5).Value = "2021-09-13T23:49:15.416539Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:10.852781"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3286/?format=json"

